# Auto-generated edits applying the diff to Maduin_Profits workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19 (hunk 0)
$ws.Range("H19").Value = 1356.3636
$ws.Range("J19").Value = 1914.4
$ws.Range("L19").Value = 1914.4
$ws.Range("N19").Value = -2264.4

# Row 33 (hunk 1)
$ws.Range("H33").Value = 370.6111
$ws.Range("I33").Value = 141.9375
$ws.Range("J33").Value = 2200
$ws.Range("K33").Value = 141.9375
$ws.Range("L33").Value = 2200
$ws.Range("M33").Value = 87.0625
$ws.Range("N33").Value = -2658

# Row 40 (hunk 2)
$ws.Range("H40").Value = 2915
$ws.Range("I40").Value = 2499.3333
$ws.Range("J40").Value = 3330.6667
$ws.Range("K40").Value = 2499.3333
$ws.Range("L40").Value = 3330.6667
$ws.Range("M40").Value = -2324.3333
$ws.Range("N40").Value = -3680.6667

# Row 100 (hunk 3)
$ws.Range("H100").Value = 2000
$ws.Range("I100").Value = 2000
$ws.Range("K100").Value = 2000
$ws.Range("M100").Value = -1459

# Row 118 (hunk 4)
$ws.Range("J118").Value = 1000
$ws.Range("L118").Value = 3000
$ws.Range("N118").Value = -6314

# Row 137 (hunk 5)
$ws.Range("H137").Value = 2367.6667
$ws.Range("I137").Value = 1901.5
$ws.Range("J137").Value = 3300
$ws.Range("K137").Value = 5704.5
$ws.Range("L137").Value = 9900
$ws.Range("M137").Value = -3154.5
$ws.Range("N137").Value = -15000

$ws = $wb.Worksheets.Item("ARM")
# Row 6 (hunk 6)
$ws.Range("H6").Value = 800
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 800
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 800
$ws.Range("M6").ClearContents()
$ws.Range("N6").Value = -1146

# Row 61 (hunk 7)
$ws.Range("H61").Value = 3598.8
$ws.Range("I61").Value = 1497.5
$ws.Range("K61").Value = 1497.5
$ws.Range("M61").Value = -1285.5

# Row 74 (hunk 8)
$ws.Range("H74").Value = 838.5333000000001
$ws.Range("I74").Value = 823.25
$ws.Range("J74").Value = 899.6667
$ws.Range("K74").Value = 823.25
$ws.Range("L74").Value = 899.6667
$ws.Range("M74").Value = 50.75
$ws.Range("N74").Value = -2647.6667

# Row 77 (hunk 9)
$ws.Range("H77").Value = 838.5333000000001
$ws.Range("I77").Value = 823.25
$ws.Range("J77").Value = 899.6667
$ws.Range("K77").Value = 4116.25
$ws.Range("L77").Value = 4498.3335
$ws.Range("M77").Value = 251.75
$ws.Range("N77").Value = -13234.3335

# Row 96 (hunk 10)
$ws.Range("H96").Value = 20014880
$ws.Range("J96").Value = 20014880
$ws.Range("L96").Value = 20014880
$ws.Range("N96").Value = -20020372

# Row 110 (hunk 11)
$ws.Range("H110").Value = 2862
$ws.Range("I110").Value = 3053.6667
$ws.Range("K110").Value = 3053.6667
$ws.Range("M110").Value = -1008.6667

# Row 132 (hunk 12)
$ws.Range("H132").Value = 1669.2
$ws.Range("I132").Value = 1524.125
$ws.Range("K132").Value = 4572.375
$ws.Range("M132").Value = -2042.375

# Row 136 (hunk 13)
$ws.Range("H136").Value = 3598.8
$ws.Range("I136").Value = 1497.5
$ws.Range("K136").Value = 4492.5
$ws.Range("M136").Value = -1942.5

$ws = $wb.Worksheets.Item("BSM")
# Row 30 (hunk 14)
$ws.Range("H30").Value = 1799
$ws.Range("J30").Value = 1799
$ws.Range("L30").Value = 1799
$ws.Range("N30").Value = -2049

# Row 94 (hunk 15)
$ws.Range("H94").Value = 979.5238000000001
$ws.Range("I94").Value = 768.94116
$ws.Range("K94").Value = 768.94116
$ws.Range("M94").Value = -317.94116

# Row 134 (hunk 16)
$ws.Range("H134").Value = 909.26666
$ws.Range("I134").Value = 909.26666
$ws.Range("K134").Value = 2727.79998
$ws.Range("M134").Value = -192.7999799999998

$ws = $wb.Worksheets.Item("CRP")
# Row 12 (hunk 17)
$ws.Range("H12").Value = 15736.728
$ws.Range("I12").Value = 3221.8
$ws.Range("J12").Value = 26165.834
$ws.Range("K12").Value = 3221.8
$ws.Range("L12").Value = 26165.834
$ws.Range("M12").Value = -3051.8
$ws.Range("N12").Value = -26505.834

# Row 31 (hunk 18)
$ws.Range("H31").Value = 2779.9167
$ws.Range("I31").Value = 2506.5557
$ws.Range("K31").Value = 2506.5557
$ws.Range("M31").Value = -2211.5557

# Row 34 (hunk 19)
$ws.Range("H34").Value = 2779.9167
$ws.Range("I34").Value = 2506.5557
$ws.Range("K34").Value = 2506.5557
$ws.Range("M34").Value = -2304.5557

# Row 92 (hunk 20)
$ws.Range("H92").Value = 44246.75
$ws.Range("J92").Value = 44246.75
$ws.Range("L92").Value = 44246.75
$ws.Range("N92").Value = -49238.75

# Row 96 (hunk 21)
$ws.Range("H96").Value = 19379.8
$ws.Range("J96").Value = 19379.8
$ws.Range("L96").Value = 19379.8
$ws.Range("N96").Value = -24871.8

# Row 99 (hunk 22)
$ws.Range("H99").Value = 2990
$ws.Range("I99").Value = 2879
$ws.Range("K99").Value = 2879
$ws.Range("M99").Value = -1381

# Row 122 (hunk 23)
$ws.Range("H122").Value = 913.7143
$ws.Range("J122").Value = 1264
$ws.Range("L122").Value = 3792
$ws.Range("N122").Value = -8692

# Row 126 (hunk 24)
$ws.Range("H126").Value = 2990
$ws.Range("I126").Value = 2879
$ws.Range("K126").Value = 8637
$ws.Range("M126").Value = -6167

$ws = $wb.Worksheets.Item("CUL")
# Row 2 (hunk 25)
$ws.Range("H2").Value = 60.8
$ws.Range("J2").Value = 126.85714
$ws.Range("L2").Value = 761.14284
$ws.Range("N2").Value = -987.14284

# Row 34 (hunk 26)
$ws.Range("H34").Value = 1330.5385
$ws.Range("J34").Value = 1799.6666
$ws.Range("L34").Value = 5398.9998
$ws.Range("N34").Value = -5566.9998

# Row 39 (hunk 27)
$ws.Range("H39").Value = 3000
$ws.Range("I39").Value = 3000
$ws.Range("K39").Value = 9000
$ws.Range("M39").Value = -8706

# Row 44 (hunk 28)
$ws.Range("H44").Value = 2927.3333
$ws.Range("I44").Value = 461.5
$ws.Range("J44").Value = 4900
$ws.Range("K44").Value = 1384.5
$ws.Range("L44").Value = 14700
$ws.Range("M44").Value = -986.5
$ws.Range("N44").Value = -15496

# Row 55 (hunk 29)
$ws.Range("H55").Value = 2550
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()

# Row 68 (hunk 30)
$ws.Range("H68").Value = 2220.5386
$ws.Range("I68").Value = 2181
$ws.Range("J68").Value = 2283.8
$ws.Range("K68").Value = 6543
$ws.Range("L68").Value = 6851.400000000001
$ws.Range("M68").Value = -5732
$ws.Range("N68").Value = -8473.400000000001

# Row 71 (hunk 31)
$ws.Range("H71").Value = 2220.5386
$ws.Range("I71").Value = 2181
$ws.Range("J71").Value = 2283.8
$ws.Range("K71").Value = 19629
$ws.Range("L71").Value = 20554.2
$ws.Range("M71").Value = -15573
$ws.Range("N71").Value = -28666.2

$ws = $wb.Worksheets.Item("GSM")
# Row 11 (hunk 32)
$ws.Range("H11").Value = 938.8
$ws.Range("I11").Value = 1000
$ws.Range("J11").Value = 923.5
$ws.Range("K11").Value = 1000
$ws.Range("L11").Value = 923.5
$ws.Range("M11").Value = -861
$ws.Range("N11").Value = -1201.5

# Row 14 (hunk 33)
$ws.Range("H14").Value = 146459.88
$ws.Range("I14").Value = 500450
$ws.Range("J14").Value = 28463.166
$ws.Range("K14").Value = 500450
$ws.Range("L14").Value = 28463.166
$ws.Range("M14").Value = -500282
$ws.Range("N14").Value = -28799.166

# Row 122 (hunk 34)
$ws.Range("H122").Value = 4170.7144
$ws.Range("I122").Value = 4032.5
$ws.Range("K122").Value = 12097.5
$ws.Range("M122").Value = -9647.5

# Row 126 (hunk 35)
$ws.Range("H126").Value = 6999.2
$ws.Range("I126").Value = 6249
$ws.Range("K126").Value = 18747
$ws.Range("M126").Value = -16277

$ws = $wb.Worksheets.Item("LTW")
# Row 19 (hunk 36)
$ws.Range("H19").Value = 1100
$ws.Range("I19").Value = 1000
$ws.Range("J19").Value = 1200
$ws.Range("K19").Value = 1000
$ws.Range("L19").Value = 1200
$ws.Range("M19").Value = -830
$ws.Range("N19").Value = -1540

# Row 40 (hunk 37)
$ws.Range("H40").Value = 4095
$ws.Range("I40").Value = 4095
$ws.Range("K40").Value = 4095
$ws.Range("M40").Value = -3959

# Row 136 (hunk 38)
$ws.Range("H136").Value = 5601.3335
$ws.Range("I136").Value = 5921.75
$ws.Range("J136").Value = 5235.143
$ws.Range("K136").Value = 17765.25
$ws.Range("L136").Value = 15705.429
$ws.Range("M136").Value = -15215.25
$ws.Range("N136").Value = -20805.429

$ws = $wb.Worksheets.Item("WVR")
# Row 29 (hunk 39)
$ws.Range("H29").Value = 12019.8
$ws.Range("I29").Value = 12999
$ws.Range("J29").Value = 11775
$ws.Range("K29").Value = 12999
$ws.Range("L29").Value = 11775
$ws.Range("M29").Value = -12709
$ws.Range("N29").Value = -12355

# Row 92 (hunk 40)
$ws.Range("H92").Value = 29800
$ws.Range("J92").Value = 29800
$ws.Range("L92").Value = 29800
$ws.Range("N92").Value = -34792

# Row 132 (hunk 41)
$ws.Range("H132").Value = 1906
$ws.Range("I132").Value = 1749.8
$ws.Range("J132").Value = 2166.3333
$ws.Range("K132").Value = 5249.4
$ws.Range("L132").Value = 6498.999899999999
$ws.Range("M132").Value = -2719.4
$ws.Range("N132").Value = -11558.9999

# Row 136 (hunk 42)
$ws.Range("H136").Value = 1050
$ws.Range("I136").Value = 1050
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 3150
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -600
$ws.Range("N136").ClearContents()
